# Commit: add the NA's under duplicate_image_filename
# Column E is "duplicate_image_filename" (see header row 1).
# Fill rows 2 through 21 of column E with the text "NA".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2:E21").Value = "NA"
